# Issue #57: Make genre required with PBCore controlled vocabulary.
#
# Add a new "Genre" column to the batch-ingest manifest fixture so the
# spreadsheet exercises the (now required) genre override used during
# bibliographic import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("Q2").Value = "Genre"

# Sample genre values (PBCore controlled vocabulary) for the two
# existing data rows.
$ws.Range("Q3").Value = "Aviation"
$ws.Range("Q4").Value = "Travel"

# Leave the selection on the newly-added cell, matching where the
# editor's cursor ended up after adding the column.
$ws.Range("Q4").Select() | Out-Null
